$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 446 (pushes the old rows 446-468 down to 447-469),
# matching the alphabetical sort position of the new dataset entry.
$ws.Rows(446).Insert()

# Fill in the new dataset row: "US global military deployments"
$ws.Range("A446").Value2 = "US global military deployments"
$ws.Range("B446").Value2 = "international relations"
$ws.Range("C446").Value2 = "https://doi.org/10.1177%2F07388942211030885"
$ws.Range("F446").Value2 = 1
$ws.Range("G446").Value2 = 1
$ws.Range("H446").Value2 = 1
$ws.Range("I446").Value2 = 1
$ws.Range("J446").Value2 = 1
$ws.Range("K446").Value2 = 1950
$ws.Range("L446").Value2 = 2020
$ws.Range("M446").Value2 = "online"
$ws.Range("N446").Value2 = "no"
$ws.Range("O446").Value2 = 1
$ws.Range("W446").Value2 = "countryname"
$ws.Range("X446").Value2 = "year"
$ws.Range("Y446").Value2 = "ccode"
$ws.Range("Z446").Value2 = "10.1177%2F07388942211030885"
$ws.Range("D446").Value2 = "alliances, overseas deployments, troop deployments, US foreign policy"
$ws.Range("AB446").Value2 = 20210808

# Give the link cell the same hyperlink style as the rest of column C
$ws.Hyperlinks.Add($ws.Range("C446"), "https://doi.org/10.1177%2F07388942211030885")
$ws.Range("C447").Copy()
$ws.Range("C446").PasteSpecial(-4122)

# Restore the selection the author left the sheet with
$ws.Range("A2:AC469").Select()
